$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values per repulled data / mean calculation
$updates = @{
    2  = -3
    3  = 3
    4  = -3
    5  = -1
    6  = 7
    8  = 4
    9  = 1
    10 = 3
    11 = 3
    12 = -2
    13 = 4
    14 = 2
    15 = 1
    16 = 1
    18 = 4
    19 = 2
    20 = 2
    21 = 2
    25 = -5
    26 = -1
    27 = 2
    28 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
